# Weekly data refresh: a new observation is prepended to the price history
# for row 50 (the Ciboulette series keeps its most-recent-first ordering),
# pushing all subsequent rows (50-309) down by one and carrying the former
# last row (309) to the new bottom row (310).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 50..309 down to 51..310, inserting a blank row at 50.
$ws.Rows(50).Insert()

# The new row 50 starts life as a duplicate of the row now below it (the
# old row 50, shifted to 51) so every column keeps its prior value...
$ws.Range("A51:R51").Copy()
$ws.Range("A50").PasteSpecial()

# ...except the date, which is the new weekly data point.
$ws.Range("D50").Value = 44687
